$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M1").Value = "age_80"
$ws.Range("L1").Value = "age_70"
$ws.Range("K1").Value = "age_60"
$ws.Range("J1").Value = "age_50"
$ws.Range("I1").Value = "age_40"
$ws.Range("H1").Value = "age_30"
$ws.Range("G1").Value = "age_20"
$ws.Range("F1").Value = "age_10"
$ws.Range("E1").Value = "age_0"
$ws.Range("D1").Value = "female"
$ws.Range("C1").Value = "male"
$ws.Range("B1").Value = "total"
$ws.Range("A1").Value = "date"
$ws.Range("A2").Select() | Out-Null
